$d = $word.ActiveDocument
$r = $d.Range(0, $d.Content.End)
$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>-- File: progymdb.sql</w:t></w:r></w:p><w:p><w:r><w:t>-- Path: ./progymdb.sql</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- CREATE DATABASE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Create the main database for the ProGym application if it doesn't already exist,</w:t></w:r></w:p><w:p><w:r><w:t>-- using UTF-8 character set for full Unicode support.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE DATABASE IF NOT EXISTS progymdb</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    CHARACTER SET utf8mb4</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    COLLATE utf8mb4_unicode_ci;</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- Switch the connection to the ProGym database for subsequent operations</w:t></w:r></w:p><w:p><w:r><w:t>USE progymdb;</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- USERS TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Stores each user's account details and profile information.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS users (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,         -- Unique user identifier</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    username VARCHAR(255) NOT NULL UNIQUE,    -- Login name, must be unique</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    email VARCHAR(255) NOT NULL UNIQUE,       -- User email, must be unique</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    password VARCHAR(255) NOT NULL,           -- Hashed password for authentication</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    role ENUM('user', 'premium') NOT NULL DEFAULT 'user',  -- Access level</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    height_cm INT,                            -- User height in centimeters</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    weight_kg INT,                            -- User weight in kilograms</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    age INT,                                  -- User age in years</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    gender ENUM('male', 'female', 'other'),   -- User gender</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    goal ENUM('cutting', 'bulking', 'maintenance'),  -- Fitness goal</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    created_at TIMESTAMP DEFAULT CURRENT_TIMESTAMP  -- Record creation timestamp</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- WORKOUTS TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Logs individual workout entries for progress tracking.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS workouts (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,        -- Unique workout entry ID</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    user_id INT NOT NULL,                    -- Which user performed the workout</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    category VARCHAR(50) NOT NULL,           -- Workout category (e.g., chest, legs)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    exercise VARCHAR(100) NOT NULL,          -- Name of the exercise performed</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    sets INT NOT NULL,                       -- Number of sets completed</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    reps INT NOT NULL,                       -- Number of repetitions per set</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    weight DECIMAL(5,2) NOT NULL,            -- Weight used, in kilograms</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    date DATE NOT NULL,                      -- Date of the workout</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    created_at TIMESTAMP DEFAULT CURRENT_TIMESTAMP,  -- Entry creation timestamp</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    FOREIGN KEY (user_id) REFERENCES users(id) ON DELETE CASCADE  -- Link to users table</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- MEALS TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Records daily meals for calorie and macro nutrient tracking.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS meals (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,        -- Unique meal entry ID</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    user_id INT NOT NULL,                    -- Which user logged the meal</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    date DATE NOT NULL,                      -- Date of the meal</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    calories INT,                            -- Total calories consumed</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    protein INT,                             -- Protein in grams</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    carbs INT,                               -- Carbohydrates in grams</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    fat INT,                                 -- Fat in grams</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    meal_type ENUM('breakfast', 'lunch', 'dinner') DEFAULT 'lunch',  -- Meal category</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    created_at TIMESTAMP DEFAULT CURRENT_TIMESTAMP,  -- Entry creation timestamp</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    FOREIGN KEY (user_id) REFERENCES users(id) ON DELETE CASCADE  -- Link to users table</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- WATER INTAKE TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Tracks daily hydration levels by total water consumed.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS water_intake (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,        -- Unique water intake entry ID</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    user_id INT NOT NULL,                    -- Which user logged the intake</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    date DATE NOT NULL,                      -- Date of water consumption</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    total_ml INT NOT NULL,                   -- Total water consumed in milliliters</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    created_at TIMESTAMP DEFAULT CURRENT_TIMESTAMP,  -- Entry creation timestamp</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    FOREIGN KEY (user_id) REFERENCES users(id) ON DELETE CASCADE  -- Link to users table</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- FOODS TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Master list of foods for meal planning and calorie lookup.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS foods (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,        -- Unique food item ID</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    name VARCHAR(100) NOT NULL,              -- Common name of the food</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    calories FLOAT NOT NULL,                 -- Calories per 100g</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    protein_per_100g FLOAT NOT NULL,         -- Protein per 100g</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    carbs_per_100g FLOAT NOT NULL,           -- Carbs per 100g</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    fat_per_100g FLOAT NOT NULL              -- Fat per 100g</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- MEAL PLANS TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Stores metadata for autogenerated meal plans.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS meal_plans (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,        -- Unique meal plan ID</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    user_id INT NOT NULL,                    -- Owner of the meal plan</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    meals_per_day INT NOT NULL,              -- Number of meals planned per day</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    total_calories INT,                      -- Total daily calories in the plan</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    total_protein FLOAT,                     -- Total daily protein</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    total_carbs FLOAT,                       -- Total daily carbs</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    total_fat FLOAT,                         -- Total daily fat</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    goal ENUM('cutting', 'bulking', 'maintenance'),  -- User's fitness goal for the plan</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    date DATE NOT NULL,                      -- Date the plan was generated</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    created_at TIMESTAMP DEFAULT CURRENT_TIMESTAMP,  -- Creation timestamp</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    FOREIGN KEY (user_id) REFERENCES users(id) ON DELETE CASCADE  -- Link to users table</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- MEAL PLAN ITEMS TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Details which foods and quantities make up each meal in a plan.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS meal_plan_items (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,        -- Unique item entry ID</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    meal_plan_id INT NOT NULL,               -- Which meal plan this item belongs to</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    meal_number INT NOT NULL,                -- Sequence number of the meal (e.g., 1 for breakfast)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    food_id INT NOT NULL,                    -- Which food item is used</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    quantity_grams INT NOT NULL,             -- Quantity of food in grams</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    created_at TIMESTAMP DEFAULT CURRENT_TIMESTAMP,  -- Entry creation timestamp</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    FOREIGN KEY (meal_plan_id) REFERENCES meal_plans(id) ON DELETE CASCADE,  -- Link to meal_plans</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    FOREIGN KEY (food_id) REFERENCES foods(id) ON DELETE CASCADE  -- Link to foods table</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- GYMBOT SESSIONS TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Tracks individual chat sessions with the GymBot feature.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS gymbot_sessions (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,        -- Unique session ID</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    user_id INT NOT NULL,                    -- Which user started the session</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    title VARCHAR(255) DEFAULT 'Untitled Chat',  -- User-assigned session title</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    created_at TIMESTAMP DEFAULT CURRENT_TIMESTAMP,  -- Session start timestamp</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    FOREIGN KEY (user_id) REFERENCES users(id) ON DELETE CASCADE  -- Link to users table</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- GYMBOT MESSAGES TABLE</w:t></w:r></w:p><w:p><w:r><w:t>-- ===========================</w:t></w:r></w:p><w:p><w:r><w:t>-- Stores each message exchanged between user and GymBot within a session.</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE IF NOT EXISTS gymbot_messages (</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    id INT AUTO_INCREMENT PRIMARY KEY,        -- Unique message ID</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    session_id INT NOT NULL,                 -- Which session this message belongs to</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    sender ENUM('user', 'bot') NOT NULL,     -- Who sent the message</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    text TEXT NOT NULL,                      -- Content of the message</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    created_at TIMESTAMP DEFAULT CURRENT_TIMESTAMP,  -- Message timestamp</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    FOREIGN KEY (session_id) REFERENCES gymbot_sessions(id) ON DELETE CASCADE  -- Link to sessions</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($newBodyXml)
Write-Host "Paragraphs after:" $d.Paragraphs.Count
